$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.610.57"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.112.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +1.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.82"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5244"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4503"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.41"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08986"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.168"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.36"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.117.73"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.774"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.052"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.84"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001159"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.013"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06697"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.31"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.012"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.300"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.674.42"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.76"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.387"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.368.33"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.26"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.27"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.534"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.72"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.189"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1070"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.637"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.350"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.942"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.858"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02645"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.57"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6862"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.257"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.93"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6405"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.305"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000368"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +10.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.700"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.65"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.72%  "
